$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "CVS Organic Multivitamin - (Product B - Yellow Gummy)"
$ws.Range("A12").Value = "Spring Valley Adult Organic Multivitamin (Product A- Red Gummy)"
$ws.Range("B11").Value = "CVS Health Organic Womens Multi"
$ws.Range("B12").Value = "Spring Valley Womens Multi"

$ws.PageSetup.Orientation = 1
$null = $ws.Range("E13").Select()
